$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update 想去人数 (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 684
$ws1.Range("F4").Value = 230
$ws1.Range("F5").Value = 1958
$ws1.Range("F6").Value = 44
$ws1.Range("F7").Value = 3328
$ws1.Range("F9").Value = 801

# Sheet "全部类型" (All Types) - update 想去人数 (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 684
$ws4.Range("F5").Value = 230
$ws4.Range("F6").Value = 1958
$ws4.Range("F7").Value = 44
$ws4.Range("F8").Value = 3328
$ws4.Range("F10").Value = 801
